# Updates cryptos list data (Coin/Link/Price/Volume(1h) columns) to match the
# latest scrape. Price values that look like plain numbers are written with a
# leading apostrophe (quote-prefix) so Excel keeps them as text, exactly like
# the original cells (e.g. "586.06" must stay text, not become the number 586.06).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.379.95'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '3.317.37'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''586.06'
$ws.Range("E5").Value = '  +2.16%  '
$ws.Range("D6").Value = '''182.46'
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("D7").Value = '''0.644'
$ws.Range("E7").Value = '  +7.68%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -2.23%  '
$ws.Range("D10").Value = '''6.79'
$ws.Range("E10").Value = '  +2.31%  '
$ws.Range("D11").Value = '''0.403'
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").Value = '3.893.88'
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("E13").Value = '  -4.42%  '
$ws.Range("D14").Value = '66.408.50'
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("D15").Value = '''26.41'
$ws.Range("E15").Value = '  -2.93%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '''0.0000164'
$ws.Range("E16").Value = '  -2.20%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.313.66'
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").Value = '''430.53'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("D19").Value = '''13.32'
$ws.Range("E19").Value = '  -2.52%  '
$ws.Range("E20").Value = '  -2.67%  '
$ws.Range("E21").Value = '  -2.62%  '
$ws.Range("D22").Value = '''72.31'
$ws.Range("E22").Value = '  -1.73%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = '''5.70'
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("D25").Value = '3.439.04'
$ws.Range("E25").Value = '  -0.75%  '
$ws.Range("E26").Value = '  -0.99%  '
$ws.Range("E27").Value = '  +2.66%  '
$ws.Range("E28").Value = '  -3.80%  '
$ws.Range("D29").Value = '''9.04'
$ws.Range("E29").Value = '  -0.82%  '
$ws.Range("D30").Value = '''1.00'
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("E32").Value = '  -1.77%  '
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("D34").Value = '''5.23'
$ws.Range("E34").Value = '  -1.92%  '
$ws.Range("D35").Value = '''6.64'
$ws.Range("E35").Value = '  -2.84%  '
$ws.Range("D36").Value = '''1.20'
$ws.Range("E36").Value = '  -3.10%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '''159.11'
$ws.Range("E37").Value = '  -0.50%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '''1.47'
$ws.Range("E38").Value = '  -2.78%  '
$ws.Range("E39").Value = '  -1.50%  '
$ws.Range("D40").Value = '''26.85'
$ws.Range("E40").Value = '  -1.54%  '
$ws.Range("D41").Value = '2.877.81'
$ws.Range("E41").Value = '  +1.61%  '
$ws.Range("D42").Value = '''0.771'
$ws.Range("E42").Value = '  -2.61%  '
$ws.Range("E43").Value = '  -2.27%  '
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("E45").Value = '  -1.38%  '
$ws.Range("D46").Value = '''6.04'
$ws.Range("E46").Value = '  -2.64%  '
$ws.Range("E47").Value = '  -2.00%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").Value = '''23.43'
$ws.Range("E48").Value = '  -4.35%  '
$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D49").Value = '''316.84'
$ws.Range("E49").Value = '  -2.42%  '
$ws.Range("E50").Value = '  -0.72%  '
$ws.Range("E51").Value = '  +3.70%  '

$wb.Save()
